# "Generate Report for Handoff" - the handoff just completed, so every
# "In Translation" status flips to "Ready for handoff" and the associated
# handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-locale status + latest handoff-xliff-generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 16:45:56"

# zh-cn detail sheet: status + latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 16:45:51"

# de-de detail sheet: status + latest handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 16:45:56"

# The "Status" columns grew wider because "Ready for handoff" is longer
# than "In Translation" - mirror Excel's auto-sizing of those columns.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
